# Generate Report for Handoff
# Replaces the stale 04678046.../44e5c74c... handoff-cycle data with the
# new 38042f86.../ffffba476e03... cycle, flips Status to "Ready for
# handoff", and drops the now-empty duplicate "Latest Target
# File"/"Latest Handback File" columns (F/G) from the per-locale sheets.

$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1) Update hyperlink display text in-place (match by the hyperlink's
#    own Range address so we don't disturb ones we're not touching).
# ---------------------------------------------------------------------
function Set-HyperlinkText($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

$newMdA = "38042f86-8738-4748-8ad0-7599bab8d2e5.md"
$newMdB = "ffffba476e03-40ca-4322-9c2c-80847e2c1d90.md"
$newXlfZh = "38042f86-8738-4748-8ad0-7599bab8d2e5.e5d31397643a3be6beeb356f600c1a98f21eadc0.zh-cn.xlf"
$newXlfDe = "38042f86-8738-4748-8ad0-7599bab8d2e5.e5d31397643a3be6beeb356f600c1a98f21eadc0.de-de.xlf"

Set-HyperlinkText $ovw '$A$2' $newMdA
Set-HyperlinkText $ovw '$A$3' $newMdB

Set-HyperlinkText $zh '$A$2' $newMdA
Set-HyperlinkText $zh '$D$2' $newXlfZh
Set-HyperlinkText $zh '$A$3' $newMdB
Set-HyperlinkText $zh '$D$3' $newXlfZh

Set-HyperlinkText $de '$A$2' $newMdA
Set-HyperlinkText $de '$D$2' $newXlfDe
Set-HyperlinkText $de '$A$3' $newMdB
Set-HyperlinkText $de '$D$3' $newXlfDe

# ---------------------------------------------------------------------
# 2) Remove the hyperlinks anchored to the columns we're about to blank
#    out (F/G on both locale sheets). Delete one at a time and re-query
#    the live collection each time -- deleting mid-foreach reindexes it.
# ---------------------------------------------------------------------
function Remove-HyperlinkAt($ws, $addr) {
    $again = $true
    while ($again) {
        $again = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $addr) {
                $hl.Delete()
                $again = $true
                break
            }
        }
    }
}

Remove-HyperlinkAt $zh '$F$2'
Remove-HyperlinkAt $zh '$G$2'
Remove-HyperlinkAt $zh '$F$3'
Remove-HyperlinkAt $zh '$G$3'

Remove-HyperlinkAt $de '$F$2'
Remove-HyperlinkAt $de '$G$2'
Remove-HyperlinkAt $de '$F$3'
Remove-HyperlinkAt $de '$G$3'

# ---------------------------------------------------------------------
# 3) Blank the now-orphaned F/G cells entirely (not just ClearContents,
#    which would leave an empty styled <c> behind).
# ---------------------------------------------------------------------
$zh.Range("F2:G2").Clear()
$zh.Range("F3:G3").Clear()
$de.Range("F2:G2").Clear()
$de.Range("F3:G3").Clear()

# ---------------------------------------------------------------------
# 4) Update the cell text itself (the diff's shared-string churn is a
#    side effect of these values changing).
# ---------------------------------------------------------------------
$status = "Ready for handoff"
$handoffDate = "2016-03-22 19:07:20"
$handoffDatetimeZh = "2016-03-22 19:07:16"
$handbackReset = "0001-01-01 00:00:00"

# Overview
$ovw.Range("A2").Value = $newMdA
$ovw.Range("B2").Value = $status
$ovw.Range("C2").Value = $status
$ovw.Range("D2").Value = $handoffDate
$ovw.Range("A3").Value = $newMdB
$ovw.Range("B3").Value = $status
$ovw.Range("C3").Value = $status
$ovw.Range("D3").Value = $handoffDate

# zh-cn
$zh.Range("A2").Value = $newMdA
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $newXlfZh
$zh.Range("E2").Value = $handoffDatetimeZh
$zh.Range("H2").Value = $handbackReset
$zh.Range("A3").Value = $newMdB
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $newXlfZh
$zh.Range("E3").Value = $handoffDatetimeZh
$zh.Range("H3").Value = $handbackReset

# de-de
$de.Range("A2").Value = $newMdA
$de.Range("C2").Value = $status
$de.Range("D2").Value = $newXlfDe
$de.Range("E2").Value = $handoffDate
$de.Range("H2").Value = $handbackReset
$de.Range("A3").Value = $newMdB
$de.Range("C3").Value = $status
$de.Range("D3").Value = $newXlfDe
$de.Range("E3").Value = $handoffDate
$de.Range("H3").Value = $handbackReset

Write-Host "Applied handoff-report regeneration."
